$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Insert two new blank rows at row 10 (pushes the old "Option" /
#    "Override" blocks that started at row 12 down to row 14, and
#    gives us two fresh rows to hold the TimeTrial entries that used
#    to live in rows 6 and 7 - ProductionNumber 9001 / 9002).
# ------------------------------------------------------------------
$ws.Rows.Item(10).Resize(2).Insert()

# ------------------------------------------------------------------
# 2. Re-create the old row 6 / row 7 TimeTrial data (ProductionNumber
#    9001 / 9002) in the newly inserted rows 10 / 11.
# ------------------------------------------------------------------
$ws.Range("G10").Value = 9001
$ws.Range("H10").Value = 151000
$ws.Range("I10").Value = [DateTime]"2018-06-04"
$ws.Range("J10").Value = 3.85
$ws.Range("K10").Value = 2.1
$ws.Range("L10").Value = 1.75
$ws.Range("M10").Value = "A1C1A022"
$ws.Range("N10").Value = 0

$ws.Range("G11").Value = 9002
$ws.Range("H11").Value = 151000
$ws.Range("I11").Value = [DateTime]"2018-06-05"
$ws.Range("J11").Value = 3.6
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 1.5
$ws.Range("M11").Value = "A1C1A022"
$ws.Range("N11").Value = 0

# ------------------------------------------------------------------
# 3. Remove the TotalTime (column J) formulas on rows 3-5 and freeze
#    them as plain numbers (values are unchanged except J4).
# ------------------------------------------------------------------
$ws.Range("J3").Value = 6.35
$ws.Range("J4").Value = 6.25
$ws.Range("J5").Value = 2.9

# ------------------------------------------------------------------
# 4. Row 6 becomes a brand new trial (ProductionNumber 8004).
# ------------------------------------------------------------------
$ws.Range("G6").Value = 8004
$ws.Range("H6").Value = 150001
$ws.Range("I6").Value = [DateTime]"2018-06-06"
$ws.Range("J6").Value = 3.97
$ws.Range("K6").Value = 1.87
$ws.Range("L6").Value = 0.95
$ws.Range("M6").Value = "A1C1A002"
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = "PATX"

# ------------------------------------------------------------------
# 5. Row 7 becomes a brand new trial (ProductionNumber 8005).
# ------------------------------------------------------------------
$ws.Range("G7").Value = 8005
$ws.Range("H7").Value = 150001
$ws.Range("I7").Value = [DateTime]"2018-06-07"
$ws.Range("J7").Value = 5.07
$ws.Range("K7").Value = 2.0299999999999998
$ws.Range("L7").Value = 1.0900000000000001
$ws.Range("M7").Value = "A1C1A002"
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = "PATX"

# ------------------------------------------------------------------
# 6. Row 8 gains a TimeTrial entry too (ProductionNumber 8006) and
#    loses its previous custom row height.
# ------------------------------------------------------------------
$ws.Range("G8").Value = 8006
$ws.Range("H8").Value = 150001
$ws.Range("I8").Value = [DateTime]"2018-06-07"
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = "A1C1A002"
$ws.Range("N8").Value = 0
$ws.Rows.Item(8).AutoFit()

# ------------------------------------------------------------------
# 7. Row 9 gains a TimeTrial entry too (ProductionNumber 8007).
# ------------------------------------------------------------------
$ws.Range("G9").Value = 8007
$ws.Range("H9").Value = 150002
$ws.Range("I9").Value = [DateTime]"2018-06-08"
$ws.Range("J9").Value = 8.65
$ws.Range("K9").Value = 2.25
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = "A1C1A002"
$ws.Range("N9").Value = 4
$ws.Range("O9").Value = "PBCTXZ"

# ------------------------------------------------------------------
# 8. Give the "OptionsText" header (O2) the same bottom border the
#    other header rows use.
# ------------------------------------------------------------------
$ws.Range("O2").Borders.Item(9).Color = 0
$ws.Range("O2").Borders.Item(9).LineStyle = 1
